$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(13, 0.3333333333333333, 0.6510902716008173, 252),
    @(14, 0.3333333333333333, 1.417572166828469, 253),
    @(15, 0.3408662900188323, 0.4124583697988894, 254),
    @(16, 0.3333333333333333, 0.3708246387322489, 255),
    @(17, 0.3408662900188323, 3.521417708659667, 256),
    @(18, 0.3333333333333333, 2.090990100992806, 257)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
